$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2021-08-24"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 720
$ws.Range("D3").Value = "2020-12-09"
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("O3").Value = "Región de La Araucanía"
$ws.Range("P3").Value = 480
$ws.Range("D4").Value = "2021-08-20"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 18500
$ws.Range("L4").Value = 18500
$ws.Range("M4").Value = 18500
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 740
$ws.Range("D5").Value = "2021-06-15"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 20000
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 800
$ws.Range("D6").Value = "2021-08-13"
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 17000
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 680
$ws.Range("D7").Value = "2021-07-02"
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 17000
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 680
$ws.Range("D8").Value = "2020-11-24"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 11500
$ws.Range("M8").Value = 11500
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 460
$ws.Range("D9").Value = "2021-08-27"
$ws.Range("J9").Value = 170
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 19000
$ws.Range("M9").Value = 18529
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 741
$ws.Range("D11").Value = "2021-09-10"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 17000
$ws.Range("M11").Value = 17000
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 680
$ws.Range("D12").Value = "2021-09-07"
$ws.Range("J12").Value = 90
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 17000
$ws.Range("M12").Value = 17000
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 680
$ws.Range("D13").Value = "2021-07-06"
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 17000
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 680
$ws.Range("D14").Value = "2020-11-26"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 11500
$ws.Range("L14").Value = 11500
$ws.Range("M14").Value = 11500
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 460
$ws.Range("D16").Value = "2021-07-23"
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 16500
$ws.Range("L16").Value = 16500
$ws.Range("M16").Value = 16500
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 660
$ws.Range("D18").Value = "2021-09-03"
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 18000
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 720
$ws.Range("D19").Value = "2021-01-18"
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 25000
$ws.Range("O19").Value = "Región de La Araucanía"
$ws.Range("P19").Value = 1000
$ws.Range("D20").Value = "2020-12-02"
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("O20").Value = "Región de La Araucanía"
$ws.Range("P20").Value = 480
$ws.Range("D21").Value = "2021-08-10"
$ws.Range("J21").Value = 90
$ws.Range("K21").Value = 18000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 18000
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 720
$ws.Range("D22").Value = "2021-08-06"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 18000
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 720
$ws.Range("D23").Value = "2020-12-07"
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 12000
$ws.Range("O23").Value = "Región de La Araucanía"
$ws.Range("P23").Value = 480
$ws.Range("D24").Value = "2021-07-13"
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 16000
$ws.Range("M24").Value = 16000
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 640
$ws.Range("D25").Value = "2021-07-09"
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 17000
$ws.Range("M25").Value = 17000
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 680
$ws.Range("D26").Value = "2020-12-01"
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 12000
$ws.Range("O26").Value = "Región Metropolitana"
$ws.Range("P26").Value = 480
$ws.Range("D27").Value = "2021-06-22"
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 18000
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 720
$ws.Range("D28").Value = "2020-11-25"
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 11500
$ws.Range("L28").Value = 11500
$ws.Range("M28").Value = 11500
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 460
$ws.Range("D29").Value = "2021-06-04"
$ws.Range("J29").Value = 30
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = 20000
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 800
$ws.Range("D30").Value = "2021-06-18"
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 20000
$ws.Range("O30").Value = "Provincia de Limarí"
$ws.Range("P30").Value = 800
$ws.Range("D31").Value = "2021-01-19"
$ws.Range("J31").Value = 60
$ws.Range("K31").Value = 25000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 25000
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 1000
$ws.Range("D32").Value = "2021-01-07"
$ws.Range("J32").Value = 50
$ws.Range("K32").Value = 22000
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = 22000
$ws.Range("O32").Value = "Región de La Araucanía"
$ws.Range("P32").Value = 880
$ws.Range("D33").Value = "2020-11-27"
$ws.Range("J33").Value = 140
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 13000
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 520
$ws.Range("D34").Value = "2021-08-03"
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 18000
$ws.Range("M34").Value = 18000
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 720
$ws.Range("D36").Value = "2021-06-11"
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = 20000
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 800
$ws.Range("D37").Value = "2021-08-17"
$ws.Range("J37").Value = 90
$ws.Range("K37").Value = 18000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 18000
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 720
